# Updates cryptos list values (Price / Volume(1h) columns, plus the
# Maker/MXToken row swap) to match the latest scrape.
#
# Price-column (D) values that look numeric ("214.86", "0.5076", ...) get
# forced to text first -- Excel auto-converts a plain numeric-looking
# string assigned to .Value into a real number, which would both change
# the stored type and round/alter values like "25.961.62" silently.
# ClearFormats() afterwards drops the temporary "@" text format so the
# cell's style stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.961.62"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "1.640.75"
$ws.Range("E3").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5076"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2579"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06362"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.87"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07745"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.297"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").Value = "1.640.92"
$ws.Range("E13").Value = "  -0.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5469"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.22%  "

# Row 15
$ws.Range("D15").Value = "0.0₅7752"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.68%  "

# Row 17
$ws.Range("D17").Value = "25.997.92"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("E18").Value = "  -0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.463"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "196.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.966"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.149"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.82%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.891"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.95%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.69%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1271"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +11.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.867"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.38%  "

# Row 28
$ws.Range("E28").Value = "  -0.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.239"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04898"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.78%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.268"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.18%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.205"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.549"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.32%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.376"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9189"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.66%  "

# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.567"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.76%  "

# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.135.52"
$ws.Range("E37").Value = "  +0.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5543"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01571"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.12%  "

# Row 40
$ws.Range("E40").Value = "  -0.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.597"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.79%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8033"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.75"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.35%  "

# Row 44
$ws.Range("E44").Value = "  -8.06%  "

# Row 45
$ws.Range("D45").Value = "1.778.97"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4525"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.29"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.67%  "

# Row 48
$ws.Range("E48").Value = "  -0.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05190"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.523"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.91%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.19%  "
